# Rearrange / prune the columns of Sheet1:
#   old A (german)              -> new A (german)
#   old B (word_freq)           -> new D (word_freq)
#   old C (aoa_german_b_m)      -> dropped
#   old D (aoa_german_s_m)      -> dropped
#   old E (aoa_german_l_m)      -> dropped
#   old F (aoa_german_comb)     -> new C (aoa_german_comb)
#   old G (clt)                 -> new F (clt)
#   old H (english)             -> new G (english)
#   old I (aoa_mor)             -> new E (aoa_mor)
#   old J (aoa_rating_english)  -> dropped
#   old K (select)              -> new B (select)
# Resulting used range shrinks from A1:K345 to A1:G345.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 345
$oldRange = $ws.Range("A1:K" + $lastRow)
$old = $oldRange.Value2

$new = New-Object 'object[,]' $lastRow,7

for ($r = 1; $r -le $lastRow; $r++) {
    # column indices into $old (1-based): A=1 B=2 C=3 D=4 E=5 F=6 G=7 H=8 I=9 J=10 K=11
    $new[$r-1,0] = $old[$r,1]   # A <- A (german)
    $new[$r-1,1] = $old[$r,11]  # B <- K (select)
    $new[$r-1,2] = $old[$r,6]   # C <- F (aoa_german_comb)
    $new[$r-1,3] = $old[$r,2]   # D <- B (word_freq)
    $new[$r-1,4] = $old[$r,9]   # E <- I (aoa_mor)
    $new[$r-1,5] = $old[$r,7]   # F <- G (clt)
    $new[$r-1,6] = $old[$r,8]   # G <- H (english)
}

$ws.Range("A1:G" + $lastRow).Value2 = $new

# drop the now-unused old columns H:K so the used range / dimension shrinks to A:G
$ws.Range("H1:K" + $lastRow).Delete()
